{"js": "const body = context.document.body;\nconst results = body.search(\"August 11, 2022\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nif (results.items.length > 0) {\n    results.items[0].insertText(\"April 18, 2023\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"August 11, 2022\"\n$find.Replacement.Text = \"April 18, 2023\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
